# Update countries & provincias Spain
# Applies the data-refresh edit to the "Pais" worksheet:
#  - Oman's case count overtook several Balkan/Caucasus countries, so its row
#    moved up in the (descending, sorted-by-"Casos totales") table, shifting
#    Azerbaiyan / Eslovenia / Lituania / Armenia / Bosnia y Herzegovina down
#    by one row each.
#  - Several other country rows received refreshed figures.
#  - The "Datos actualizados..." timestamp string was bumped from 08:22 to 08:52.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Timestamp banner (row 1)
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 20 de Abril de 2020 a las 08:52"

# Alemania (row 8)
$ws.Cells.Item(8, 4).Value = 91500
$ws.Cells.Item(8, 5).Value = 49600

# Israel (row 25)
$ws.Cells.Item(25, 2).Value = 13654
$ws.Cells.Item(25, 3).Value = 163
$ws.Cells.Item(25, 4).Value = 3872
$ws.Cells.Item(25, 5).Value = 9609
$ws.Cells.Item(25, 6).Value = 150
$ws.Cells.Item(25, 7).Value = 1
$ws.Cells.Item(25, 8).Value = 173

# Ucrania (row 44)
$ws.Cells.Item(44, 2).Value = 5710
$ws.Cells.Item(44, 3).Value = 261
$ws.Cells.Item(44, 4).Value = 359
$ws.Cells.Item(44, 5).Value = 5200
$ws.Cells.Item(44, 7).Value = 10
$ws.Cells.Item(44, 8).Value = 151

# Kazajistan (row 67)
$ws.Cells.Item(67, 2).Value = 1757
$ws.Cells.Item(67, 3).Value = 81
$ws.Cells.Item(67, 4).Value = 413
$ws.Cells.Item(67, 5).Value = 1325

# Row 72: was Azerbaiyan, now Oman (moved up, new figures)
$ws.Cells.Item(72, 1).Value = "Oman"
$ws.Cells.Item(72, 2).Value = 1410
$ws.Cells.Item(72, 3).Value = 144
$ws.Cells.Item(72, 4).Value = 238
$ws.Cells.Item(72, 5).Value = 1165
$ws.Cells.Item(72, 6).Value = 3
$ws.Cells.Item(72, 8).Value = 7

# Row 73: was Eslovenia, now Azerbaiyan (shifted down one row)
$ws.Cells.Item(73, 1).Value = "Azerbaiyan"
$ws.Cells.Item(73, 2).Value = 1398
$ws.Cells.Item(73, 4).Value = 712
$ws.Cells.Item(73, 5).Value = 667
$ws.Cells.Item(73, 6).Value = 21
$ws.Cells.Item(73, 8).Value = 19

# Row 74: was Lituania, now Eslovenia (shifted down one row)
$ws.Cells.Item(74, 1).Value = "Eslovenia"
$ws.Cells.Item(74, 2).Value = 1330
$ws.Cells.Item(74, 3).Value = 0
$ws.Cells.Item(74, 4).Value = 192
$ws.Cells.Item(74, 5).Value = 1064
$ws.Cells.Item(74, 6).Value = 26
$ws.Cells.Item(74, 7).Value = 0
$ws.Cells.Item(74, 8).Value = 74

# Row 75: was Armenia, now Lituania (shifted down one row, with refreshed figures)
$ws.Cells.Item(75, 1).Value = "Lituania"
$ws.Cells.Item(75, 2).Value = 1326
$ws.Cells.Item(75, 3).Value = 28
$ws.Cells.Item(75, 4).Value = 242
$ws.Cells.Item(75, 5).Value = 1047
$ws.Cells.Item(75, 6).Value = 14
$ws.Cells.Item(75, 7).Value = 2
$ws.Cells.Item(75, 8).Value = 37

# Row 76: was Bosnia y Herzegovina, now Armenia (shifted down one row)
$ws.Cells.Item(76, 1).Value = "Armenia"
$ws.Cells.Item(76, 2).Value = 1291
$ws.Cells.Item(76, 4).Value = 545
$ws.Cells.Item(76, 5).Value = 726
$ws.Cells.Item(76, 6).Value = 30
$ws.Cells.Item(76, 8).Value = 20

# Row 77: was Oman, now Bosnia y Herzegovina (shifted down one row)
$ws.Cells.Item(77, 1).Value = "Bosnia y Herzegovina"
$ws.Cells.Item(77, 2).Value = 1285
$ws.Cells.Item(77, 4).Value = 347
$ws.Cells.Item(77, 5).Value = 890
$ws.Cells.Item(77, 6).Value = 4
$ws.Cells.Item(77, 8).Value = 48

# Letonia (row 90)
$ws.Cells.Item(90, 2).Value = 739
$ws.Cells.Item(90, 3).Value = 12
$ws.Cells.Item(90, 5).Value = 646

# Montenegro (row 114)
$ws.Cells.Item(114, 2).Value = 311
$ws.Cells.Item(114, 3).Value = 3
$ws.Cells.Item(114, 5).Value = 251
